$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.26"
$ws.Range("E2").Value = "'0.88%"
$ws.Range("D3").Value = "'35.86"
$ws.Range("E3").Value = "'1.47%"
$ws.Range("D4").Value = "'5.065"
$ws.Range("E4").Value = "'-0.13%"
$ws.Range("D5").Value = "'0.08065"
$ws.Range("E5").Value = "'1.75%"
$ws.Range("E6").Value = "'1.98%"
$ws.Range("D7").Value = "'4.139"
$ws.Range("E7").Value = "'2.27%"
$ws.Range("E8").Value = "'0.75%"
$ws.Range("D9").Value = "'0.9304"
$ws.Range("E9").Value = "'0.04%"
$ws.Range("D10").Value = "'0.1286"
$ws.Range("E10").Value = "'-6.72%"
$ws.Range("D11").Value = "'0.1910"
$ws.Range("E11").Value = "'0.57%"
$ws.Range("D12").Value = "'0.09237"
$ws.Range("E12").Value = "'1.12%"
$ws.Range("D13").Value = "'0.03482"
$ws.Range("E13").Value = "'1.32%"
$ws.Range("D14").Value = "'0.09877"
$ws.Range("E14").Value = "'0.44%"
$ws.Range("D15").Value = "'0.001420"
$ws.Range("E15").Value = "'1.73%"
$ws.Range("E16").Value = "'12.74%"
$ws.Range("D17").Value = "'3.609"
$ws.Range("E17").Value = "'2.19%"
$ws.Range("E18").Value = "'1.79%"
$ws.Range("E19").Value = "'-0.08%"
$ws.Range("E20").Value = "'2.45%"
$ws.Range("D21").Value = "'5.192"
$ws.Range("E21").Value = "'2.89%"
$ws.Range("E22").Value = "'5.35%"
$ws.Range("D23").Value = "'0.04422"
$ws.Range("E23").Value = "'-1.83%"
$ws.Range("E24").Value = "'1.59%"
$ws.Range("D25").Value = "'0.004726"
$ws.Range("E25").Value = "'-0.70%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'5.76%"
$ws.Range("E27").Value = "'4.21%"
$ws.Range("D39").Value = "'0.02001"
$ws.Range("E39").Value = "'7.76%"
$ws.Range("D40").Value = "'0.05130"
$ws.Range("E40").Value = "'7.65%"
$ws.Range("D41").Value = "'0.007624"
$ws.Range("E41").Value = "'3.53%"
$ws.Range("E42").Value = "'5.04%"
$ws.Range("D43").Value = "'0.1367"
$ws.Range("E43").Value = "'3.19%"
$ws.Range("D44").Value = "'0.002104"
$ws.Range("E44").Value = "'-0.41%"
$ws.Range("D45").Value = "'0.01071"
$ws.Range("E45").Value = "'-1.57%"
$ws.Range("D46").Value = "'0.00006114"
$ws.Range("E46").Value = "'-2.08%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'65.22"
$ws.Range("E48").Value = "'0.84%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.06%"
